# Generate Report for Handoff
# This script updates the localization-status workbook to reflect that the
# d7a0cf63-090d-422d-b0f4-e380aabaedff file has been re-handed-off (its
# handback was stale compared to the latest source), updating status,
# timestamps and an error detail message on all three sheets.

$wb = $excel.ActiveWorkbook

$newStatus      = "Ready for handoff"
$overviewTime   = "2016-08-16 08:48:18"
$zhHandoffTime  = "2016-08-16 08:48:13"
$deHandoffTime  = "2016-08-16 08:48:18"
$errorDetail    = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d4e4d5a5b67a5b754f17a846b482f64ba681960a/e2e/d7a0cf63-090d-422d-b0f4-e380aabaedff.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d25537b772679c9cceba755689ee079a73da84fb/e2e/d7a0cf63-090d-422d-b0f4-e380aabaedff.md."

# --- Overview sheet: row 3 is the d7a0cf63... file -----------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("G3").Value = $overviewTime

# --- zh-cn sheet: row 3 is the d7a0cf63... file ---------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("H3").Value = $zhHandoffTime
$wsZh.Range("P3").Value = $errorDetail
$wsZh.Columns.Item(16).ColumnWidth = 39.15

# --- de-de sheet: row 3 is the d7a0cf63... file ---------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("H3").Value = $deHandoffTime
$wsDe.Range("P3").Value = $errorDetail
$wsDe.Columns.Item(16).ColumnWidth = 39.15
